$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new D (price) and/or E (volume) values.
# Only cells that actually changed per the diff are listed here.
$updates = @{
    2  = @{ D = "29.417.80"; E = "  -0.87%  " }
    3  = @{ D = "1.901.57";  E = "  -0.95%  " }
    4  = @{ E = "  +0.24%  " }
    5  = @{ D = "325.87";    E = "  -2.77%  " }
    6  = @{ E = "  +0.21%  " }
    7  = @{ D = "0.4805";    E = "  +2.87%  " }
    8  = @{ D = "0.4071";    E = "  -1.09%  " }
    9  = @{ D = "0.08064";   E = "  +0.36%  " }
    10 = @{ D = "1.004";     E = "  -1.06%  " }
    11 = @{ E = "  +4.01%  " }
    12 = @{ D = "1.983.72";  E = "  +2.63%  " }
    13 = @{ D = "5.948";     E = "  -0.65%  " }
    14 = @{ D = "7.072";     E = "  -1.48%  " }
    15 = @{ D = "89.86" }
    16 = @{ E = "  +0.28%  " }
    17 = @{ D = "0.06682";   E = "  +1.33%  " }
    18 = @{ E = "  -0.36%  " }
    19 = @{ D = "17.63";     E = "  -1.18%  " }
    20 = @{ E = "  +0.22%  " }
    21 = @{ D = "29.431.49"; E = "  -0.69%  " }
    22 = @{ D = "5.535";     E = "  -0.74%  " }
    23 = @{ D = "11.77";     E = "  +0.82%  " }
    24 = @{ E = "  -2.20%  " }
    25 = @{ D = "2.141.36";  E = "  -1.06%  " }
    26 = @{ D = "155.24";    E = "  -0.33%  " }
    27 = @{ D = "19.75";     E = "  -0.58%  " }
    28 = @{ D = "6.070";     E = "  +5.63%  " }
    29 = @{ D = "2.090";     E = "  -2.49%  " }
    30 = @{ D = "118.32";    E = "  +0.67%  " }
    31 = @{ E = "  -3.49%  " }
    32 = @{ D = "0.09494";   E = "  +0.10%  " }
    33 = @{ D = "1.390";     E = "  -3.13%  " }
    34 = @{ D = "3.539";     E = "  -1.11%  " }
    35 = @{ D = "5.403";     E = "  -0.06%  " }
    36 = @{ D = "0.02252";   E = "  -0.81%  " }
    37 = @{ E = "  -1.10%  " }
    38 = @{ E = "  -0.38%  " }
    39 = @{ E = "  -0.28%  " }
    40 = @{ D = "7.892";     E = "  -6.36%  " }
    41 = @{ D = "0.1845";    E = "  -0.04%  " }
    42 = @{ D = "10.21";     E = "  -0.09%  " }
    43 = @{ D = "2.415";     E = "  +2.15%  " }
    44 = @{ D = "1.284";     E = "  +1.88%  " }
    45 = @{ D = "0.07747";   E = "  +3.13%  " }
    46 = @{ D = "12.20";     E = "  -0.54%  " }
    47 = @{ D = "0.5514";    E = "  -1.03%  " }
    48 = @{ D = "1.921";     E = "  -0.66%  " }
    49 = @{ D = "113.44";    E = "  +0.40%  " }
    50 = @{ D = "0.2945";    E = "  -2.17%  " }
    51 = @{ D = "72.21";     E = "  +1.06%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    if ($cols.ContainsKey("D")) {
        # Prices are stored as plain text (not numbers) in the source sheet,
        # e.g. "6.070" or "29.417.80". A leading apostrophe forces Excel to
        # keep the literal text (and its trailing zeros / multi-dot shape)
        # instead of silently coercing it to a Number cell. Resetting the
        # style back to Normal afterwards drops the transient "quote
        # prefix" cell style Excel applies when it sees the apostrophe, so
        # the cell's formatting stays exactly as it was before the edit.
        $cell = $ws.Range("D$row")
        $cell.Value = "'" + $cols["D"]
        $cell.Style = "Normal"
    }
    if ($cols.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.Value = "'" + $cols["E"]
        $cell.Style = "Normal"
    }
}
